# Update Daily Report: 2026-02-06
#
# Appends the 2026-02-05 (Excel serial 46058) daily snapshot to Daily_Data,
# then refreshes the two summary sheets (Today_Summary, Monthly_Stats) that
# are derived from it. Only BRINK'S, INC. actually moved today (an
# Eligible-side withdrawal of 16177.957); every other depository/region row
# simply rolls forward unchanged from the prior day (row 507).

$wb = $excel.ActiveWorkbook
$dailyData = $wb.Worksheets.Item("Daily_Data")
$todaySummary = $wb.Worksheets.Item("Today_Summary")
$monthlyStats = $wb.Worksheets.Item("Monthly_Stats")

# --- 1. Append today's rows to Daily_Data -----------------------------
# Columns: B=Region_Type, C=PREV_TOTAL, D=RECEIVED, E=WITHDRAWN,
#          F=NET_CHANGE, G=ADJUSTMENT, H=TOTAL_TODAY
$reportDate = 46058

$newRows = @(
    @{ B = "ASAHI DEPOSITORY LLC Registered"; C = 0; D = 0; E = 0; F = 0; G = 0; H = 0 }
    @{ B = "ASAHI DEPOSITORY LLC Eligible"; C = 0; D = 0; E = 0; F = 0; G = 0; H = 0 }
    @{ B = "BRINK'S, INC. Registered"; C = 76497.842; D = 0; E = 0; F = 0; G = 0; H = 76497.842 }
    @{ B = "BRINK'S, INC. Eligible"; C = 98856.745; D = 0; E = 16177.957; F = -16177.957; G = 0; H = 82678.788 }
    @{ B = "CNT DEPOSITORY, INC. Registered"; C = 1246.06; D = 0; E = 0; F = 0; G = 0; H = 1246.06 }
    @{ B = "CNT DEPOSITORY, INC. Eligible"; C = 0; D = 0; E = 0; F = 0; G = 0; H = 0 }
    @{ B = "DELAWARE DEPOSITORY Registered"; C = 1633.941; D = 0; E = 0; F = 0; G = 0; H = 1633.941 }
    @{ B = "DELAWARE DEPOSITORY Eligible"; C = 18459.584; D = 0; E = 0; F = 0; G = 0; H = 18459.584 }
    @{ B = "HSBC BANK, USA Registered"; C = 1394.758; D = 0; E = 0; F = 0; G = 0; H = 1394.758 }
    @{ B = "HSBC BANK, USA Eligible"; C = 9281.978999999999; D = 0; E = 0; F = 0; G = 0; H = 9281.978999999999 }
    @{ B = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered"; C = 2395.448; D = 0; E = 0; F = 0; G = 0; H = 2395.448 }
    @{ B = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible"; C = 0; D = 0; E = 0; F = 0; G = 0; H = 0 }
    @{ B = "JP MORGAN CHASE BANK NA Registered"; C = 114985.579; D = 0; E = 0; F = 0; G = 0; H = 114985.579 }
    @{ B = "JP MORGAN CHASE BANK NA Eligible"; C = 75484.511; D = 0; E = 0; F = 0; G = 0; H = 75484.511 }
    @{ B = "LOOMIS INTERNATIONAL (US) LLC Registered"; C = 63745.991; D = 0; E = 0; F = 0; G = 0; H = 63745.991 }
    @{ B = "LOOMIS INTERNATIONAL (US) LLC Eligible"; C = 132077.206; D = 0; E = 0; F = 0; G = 0; H = 132077.206 }
    @{ B = "MALCA-AMIT USA, LLC Registered"; C = 395.145; D = 0; E = 0; F = 0; G = 0; H = 395.145 }
    @{ B = "MALCA-AMIT USA, LLC Eligible"; C = 0; D = 0; E = 0; F = 0; G = 0; H = 0 }
    @{ B = "MANFRA, TORDELLA & BROOKES, LLC Registered"; C = 50220.42; D = 0; E = 0; F = 0; G = 0; H = 50220.42 }
    @{ B = "MANFRA, TORDELLA & BROOKES, LLC Eligible"; C = 1804.683; D = 0; E = 0; F = 0; G = 0; H = 1804.683 }
    @{ B = "STONEX PRECIOUS METALS LLC Registered"; C = 14122.765; D = 0; E = 0; F = 0; G = 0; H = 14122.765 }
    @{ B = "STONEX PRECIOUS METALS LLC Eligible"; C = 16.075; D = 0; E = 0; F = 0; G = 0; H = 16.075 }
)

$lastRow = $dailyData.UsedRange.Rows.Count
$r = $lastRow + 1
foreach ($row in $newRows) {
    $dailyData.Cells.Item($r, 1).Value = $reportDate
    $dailyData.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $dailyData.Cells.Item($r, 2).Value = $row.B
    $dailyData.Cells.Item($r, 3).Value = $row.C
    $dailyData.Cells.Item($r, 4).Value = $row.D
    $dailyData.Cells.Item($r, 5).Value = $row.E
    $dailyData.Cells.Item($r, 6).Value = $row.F
    $dailyData.Cells.Item($r, 7).Value = $row.G
    $dailyData.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# --- 2. Refresh Today_Summary for the depository that moved ----------
# BRINK'S, INC.: Eligible drops from 98856.745 to 82678.788; Registered
# is unchanged; Total_Stock is recomputed as Eligible + Registered.
$brinksEligible = 82678.788
$brinksRegistered = 76497.842
$todaySummary.Range("B3").Value = $brinksEligible
$todaySummary.Range("D3").Value = $brinksEligible + $brinksRegistered

# --- 3. Refresh Monthly_Stats (month-to-date rollups for 2026-02) ----
# Row 2: overall 2026-02 Eligible/Grand_Total rollup.
$monthlyStats.Range("B2").Value = 319802.826
$monthlyStats.Range("D2").Value = 646440.775

# Row 10: 2026-02 detail line for "BRINK'S, INC. Eligible" (WITHDRAWN/
# TOTAL_TODAY columns).
$monthlyStats.Range("D10").Value = 16177.957
$monthlyStats.Range("E10").Value = 82678.788
